$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Object_Mapping")

# Remove the "Wind_Farm" / "Wind_farm" row (row 8) entirely, shifting
# subsequent rows up by one.
$ws.Rows.Item(8).Delete()
